# Update metrics values in Sheet1 as described by the commit:
# "atualizado todo o treinamento para o novo lm" (updated all training for new lm)
#
# Columns B..Q (16 metric columns) are updated for rows 2-26.
# Rows 2-25 all share one common set of new values (they shared identical
# values before the edit too). Row 26 gets its own distinct new values.
# Note: scientific notation literals are not supported by this script
# engine, so all values are written in plain decimal form.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns in order B..Q
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")

# New values shared by rows 2 through 25
$commonValues = @{
    "B" = 0.9999824846980448
    "C" = 0.9990072066285755
    "D" = 0.9997886867889028
    "E" = 0.9999917436749221
    "F" = 0.9999561883117981
    "G" = 0.00001634976554226276
    "H" = 0.0009267290336323482
    "I" = 0.00008630143841202433
    "J" = 0.00001588459595523943
    "K" = 0.00005109301718363188
    "L" = 0.0002547310818618769
    "M" = 0.004043484331892824
    "N" = 1.000016814689877
    "O" = 0.004215623831067973
    "P" = 120.0425940012932
    "Q" = 179.767509419835
}

# New values specific to row 26
$row26Values = @{
    "B" = 0.9999824855511625
    "C" = 0.9990072063719976
    "D" = 0.9997886983169935
    "E" = 0.9999917438147951
    "F" = 0.9999561904457616
    "G" = 0.0000163489691943579
    "H" = 0.0009267292731365613
    "I" = 0.00008629673027852063
    "J" = 0.00001588432684940373
    "K" = 0.00005109052856396218
    "L" = 0.0002547238100402147
    "M" = 0.004043385857713544
    "N" = 1.000016813870884
    "O" = 0.004215521164639012
    "P" = 120.0426914176503
    "Q" = 179.7676068361921
}

for ($row = 2; $row -le 25; $row++) {
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = $commonValues[$col]
    }
}

foreach ($col in $cols) {
    $ws.Range("${col}26").Value = $row26Values[$col]
}
